$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reset password")

# Add the new row of data (row 3) - order matters for shared string table ordering
$ws.Range("B3").Value = "//p[@class='alert__des']"
$ws.Range("C3").Value = "メールアドレスは、メールアドレス形式で入力してください。"
$ws.Range("A3").Value = "InvalidEmailError"

# Update selection to A3
$ws.Range("A3").Select()
